$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sense resistor change: R mid (H8) switches from a formula (2*6.8) to a
# fixed, more reliable resistor value.
$ws.Range("H8").Value = 11

# Re-color "Andy's design values" table (G3:J12) from yellow to blue to
# flag it as updated. (0,176,240) => R + G*256 + B*65536
$ws.Range("G3:J12").Interior.Color = 15773696

# Move the active cell selection.
$ws.Range("H13").Select()
